$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.7
$ws.Range("I2").Value = 4.33
$ws.Range("J2").Value = 2.3
$ws.Range("K2").Value = 2.3
$ws.Range("L2").Value = 4.5
$ws.Range("Z2").Value = 15
$ws.Range("AE2").Value = 13
$ws.Range("AI2").Value = 23
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 9
$ws.Range("AQ2").Value = 29
$ws.Range("AU2").Value = 7.5
$ws.Range("AY2").Value = 26
$ws.Range("AZ2").Value = 67

# Row 3 updates
$ws.Range("G3").Value = 2.45
$ws.Range("I3").Value = 3.2
$ws.Range("J3").Value = 3.2
$ws.Range("K3").Value = 2
$ws.Range("N3").Value = 7.5
$ws.Range("W3").Value = 7.5
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 23
$ws.Range("AH3").Value = 9
$ws.Range("AJ3").Value = 12
$ws.Range("AO3").Value = 15
$ws.Range("AQ3").Value = 51
$ws.Range("AR3").Value = 81
$ws.Range("AW3").Value = 4.75
$ws.Range("AX3").Value = 17
$ws.Range("AZ3").Value = 51
